$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(344, 44418, 0, 4, 56.89900426742533),
    @(345, 44419, 0, 4, 56.89900426742533),
    @(346, 44420, 0, 4, 56.89900426742533),
    @(347, 44421, 1, 4, 56.89900426742533),
    @(348, 44422, 0, 4, 56.89900426742533),
    @(349, 44423, 0, 2, 28.44950213371266),
    @(350, 44424, 0, 1, 14.22475106685633),
    @(351, 44425, 1, 2, 28.44950213371266),
    @(352, 44426, 0, 2, 28.44950213371266),
    @(353, 44427, 0, 2, 28.44950213371266),
    @(354, 44428, 0, 1, 14.22475106685633),
    @(355, 44429, 0, 1, 14.22475106685633),
    @(356, 44430, 1, 2, 28.44950213371266),
    @(357, 44431, 0, 2, 28.44950213371266)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Range("A343").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
}
